# The deck's cached "datetimeFigureOut" fields (date placeholders on the
# slide master, every custom layout, and the notes master) were re-saved a
# day earlier: 2025-09-01 -> 2025-08-31. Update each cached field string to
# match (the field id/type stay the same; we're only correcting the cached
# display text, same as PowerPoint does when it re-caches the field).

$p = $ppt.ActivePresentation

# --- Slide Master: Holder 5 (idx="6", type="dt") ---
$master = $p.SlideMaster
$master.Shapes.Item(14).TextFrame.TextRange.Text = "8/31/2025"

# --- Custom Layouts (slideLayout1.xml .. slideLayout5.xml) ---
# Layout 1 "Title Slide"        -> Holder 5
$master.CustomLayouts.Item(1).Shapes.Item(4).TextFrame.TextRange.Text = "8/31/2025"
# Layout 2 "Title and Content"  -> Holder 5
$master.CustomLayouts.Item(2).Shapes.Item(4).TextFrame.TextRange.Text = "8/31/2025"
# Layout 3 "Two Content"        -> Holder 6
$master.CustomLayouts.Item(3).Shapes.Item(5).TextFrame.TextRange.Text = "8/31/2025"
# Layout 4 "Title Only"         -> Holder 4
$master.CustomLayouts.Item(4).Shapes.Item(3).TextFrame.TextRange.Text = "8/31/2025"
# Layout 5 "Blank"              -> Holder 3
$master.CustomLayouts.Item(5).Shapes.Item(2).TextFrame.TextRange.Text = "8/31/2025"

# --- Notes Master: Date Placeholder 2 (uses day-month-year format) ---
# (Updated via the HeadersFooters facade - the notes master's date
# placeholder shape does not accept direct TextFrame edits.)
$notesMaster = $p.NotesMaster
$notesMaster.HeadersFooters.DateAndTime.Text = "31-08-2025"
